# Update BOC USD rates (auto)
# Appends a new publish-time snapshot row (row 15) to the "All Published
# Values" sheet, grows the sheet's AutoFilter / _FilterDatabase range to
# cover it, and bumps the "publishes" count on the "Daily Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsAll = $wb.Worksheets.Item("All Published Values")
$wsDaily = $wb.Worksheets.Item("Daily Summary")

# --- Append the new data row -------------------------------------------
$newRow = 15
$columns = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$values = @(
    "2026-01-02",
    "2026-01-02 21:44:23",
    "697.85",
    "697.85",
    "700.79",
    "700.79",
    "702.88",
    "2026/01/02 21:44:23",
    "2026-01-02 13:53:04",
    "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"
)

for ($i = 0; $i -lt $columns.Length; $i++) {
    $cellRange = $wsAll.Range($columns[$i] + $newRow)
    # Force text storage so numeric-/date-looking strings ("697.85",
    # "2026-01-02", ...) aren't silently coerced into numbers or dates by
    # Excel's smart-entry behaviour -- the source rows store everything as
    # plain text.
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $values[$i]
    # Drop back to the default "Normal" style so the new row matches the
    # unstyled look of the existing data rows instead of keeping the
    # Text-format style we just applied.
    $cellRange.Style = "Normal"
}

# --- Re-apply the AutoFilter over the expanded range ---------------------
$wsAll.AutoFilterMode = $false
$wsAll.Range("A1:J15").AutoFilter()

# --- Expand the hidden _FilterDatabase defined name to match -------------
for ($i = 1; $i -le $wb.Names.Count(); $i++) {
    $definedName = $wb.Names.Item($i)
    if ($definedName.Name() -eq "All Published Values!_FilterDatabase") {
        $definedName.RefersTo = "='All Published Values'!`$A`$1:`$J`$15"
    }
}

# --- Update the Daily Summary "publishes" count for the day --------------
$wsDaily.Range("B4").Value = 14
